$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.130674839019775
$ws.Range("B1").Value = 1.924805521965027
$ws.Range("C1").Value = 5.501724720001221
$ws.Range("D1").Value = 0.7044932842254639
$ws.Range("E1").Value = 0.7881097197532654
